$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Second")

# New shared-string text for cell B5 (row/column info on the exposed POI "cell").
$ws.Range("B5").Value = 'This Cell is at row ${cell.rowIndex} and column ${cell.columnIndex}.'

# New shared-string text for cell D7 (wrap-text indicator), with its own
# word-wrapped style, a taller row, and a wider column.
$ws.Range("D7").Value = 'This Cell''s text is ${cell.cellStyle.wrapText ? "wrapped." : "not wrapped."}'
$ws.Range("D7").WrapText = $true

$ws.Rows.Item(7).RowHeight = 60
$ws.Columns.Item(4).ColumnWidth = 17.42
